# Refresh the scraped '想去人数' (want-to-go count) figures in column F
# across the 展览 / 演出 / 全部类型 sheets to match the latest crawl
# (gh-pages output generated at 456a3b4). Only column F values change;
# no rows/columns/sheets are added or removed and no formatting changes.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    2 = 35
    3 = 68
    4 = 26
    5 = 1153
    8 = 243
    9 = 7131
    10 = 184
    12 = 5465
    13 = 10
    15 = 0
    16 = 1095
    17 = 0
    18 = 0
    19 = 273
    21 = 0
    24 = 10112
    25 = 83
    27 = 0
    28 = 47
    29 = 0
    30 = 2106
    31 = 81
    33 = 168
    34 = 1033
    37 = 310
    38 = 0
    40 = 0
    41 = 659
    42 = 111
    43 = 168
    45 = 1076
    46 = 985
    47 = 1369
    48 = 65
    49 = 1094
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    5 = 5
    8 = 0
    13 = 0
    14 = 0
    17 = 0
    18 = 0
    21 = 0
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    2 = 0
    3 = 35
    5 = 26
    6 = 15
    7 = 1153
    9 = 147
    10 = 243
    14 = 5
    16 = 0
    17 = 72
    20 = 1095
    21 = 0
    22 = 0
    23 = 0
    24 = 206
    25 = 104
    26 = 0
    28 = 83
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    34 = 88
    35 = 1033
    36 = 0
    37 = 7
    38 = 2080
    39 = 0
    40 = 0
    42 = 659
    43 = 0
    44 = 168
    45 = 1100
    46 = 1076
    47 = 985
    48 = 0
    49 = 65
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
